# Update the "Förändrad" (Changed) date column (C) for rows 2-27
# from serial date 45212 (2023-10-13) to 45221 (2023-10-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C27").Value = 45221
